$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing entry row (row 9) to the new row (row 10)
$ws.Range("A9:C9").Copy()
$ws.Range("A10:C10").PasteSpecial(-4122)

# Fill in the new log entry values
$ws.Range("A10").Value = 43896
$ws.Range("B10").Value = 1.5
$ws.Range("C10").Value = "Sprint review du sprint 2 avec le chef de projet"

# Update selection to the next empty row as in the final workbook
$ws.Range("C11").Select()
